$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("QuantitativeMetrics")

$ws.Range("B11").Value = 0.1902584951434026
$ws.Range("B12").Value = 0.3556023046779966
$ws.Range("C12").Value = "{'codebleu': 0.3556023046779966, 'ngram_match_score': 0.1902584951434026, 'weighted_ngram_match_score': 0.2027487302131021, 'syntax_match_score': 0.5642857142857143, 'dataflow_match_score': 0.46511627906976744}"
$ws.Range("B13").Value = 0.9014018910120881
